# Refactor Code containing duplicates
# The shared string "Doe " (with trailing space) duplicated part of the
# employee's name; update the cell to use a clean "Doe" value instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Doe"

# Reflect the resulting selection on the sheet (Excel records the
# active cell/selection after an edit).
$ws.Range("B1").Select()
